# The source table (操作 / Tik接口 mapping) has a row, right after the
# "向量标准化(normalize)" row and right before the "pcl" section header row,
# whose two cells are blank placeholders. This change fills that row in with
# the Euclidean-distance (L2-Norm) mapping to vec_mul / vec_add / scalar_sqrt.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row anchored on its neighbours' text, falling back to "first
# fully blank row" if the surrounding text ever shifts.
$targetRow = $null
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $prevText = $t.Cell($r - 1, 1).Range.Text
    $curC1 = $t.Cell($r, 1).Range.Text.Trim([char]7, [char]13)
    $curC2 = $t.Cell($r, 2).Range.Text.Trim([char]7, [char]13)
    if ($prevText -like "*normalize*" -and $curC1 -eq "" -and $curC2 -eq "") {
        $targetRow = $r
        break
    }
}

if ($targetRow -eq $null) {
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        $curC1 = $t.Cell($r, 1).Range.Text.Trim([char]7, [char]13)
        $curC2 = $t.Cell($r, 2).Range.Text.Trim([char]7, [char]13)
        if ($curC1 -eq "" -and $curC2 -eq "") {
            $targetRow = $r
            break
        }
    }
}

if ($targetRow -eq $null) {
    throw "Could not locate the empty target row"
}

# --- Column 1: 欧氏距离(L2-Norm)<=>向量点乘+开平方根 ---
$cell1Range = $t.Cell($targetRow, 1).Range
$cell1Range.Collapse(1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>欧氏距离</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>L2-Norm)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>&lt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>=&gt;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>向量点乘</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>+</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>开平方根</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$cell1Range.InsertXML($xml1)

# --- Column 2: vec_mul、vec_add、scalar_sqrt ---
$cell2Range = $t.Cell($targetRow, 2).Range
$cell2Range.Collapse(1)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>v</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>ec_mul</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>、</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>v</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>ec_add</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>、</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="宋体" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>calar_sqrt</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$cell2Range.InsertXML($xml2)

Write-Host "Filled row $targetRow"

